# domibusWebAdmin Installation Guide - "working on installation guide"
#
# 1. Rename the "Required Software" heading to "Requirements".
# 2. Move the "_GoBack" bookmark (last-edit-position marker) from the end
#    of the "Upgrade from 3.1.x to 3.2" paragraph to the end of the
#    "Tomcat 7" bullet - i.e. wherever the cursor was last left by the author.

$d = $word.ActiveDocument

# --- 1. Heading text change -------------------------------------------------
$d.Content.Find.Execute("Required Software", $true, $false, $false, $false, `
                         $false, $true, 1, $false, "Requirements", 2) | Out-Null

# --- 2. Relocate the "_GoBack" bookmark -------------------------------------
# Remove the existing bookmark wherever it currently sits.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Find the end of the "Tomcat 7" run.
$target = $d.Content
$null = $target.Find.Execute("Tomcat 7", $true, $false, $false, $false, `
                              $false, $true, 1, $false, "", 0)
$target.Collapse(0)

# Inserting a genuinely empty bookmark exactly on a paragraph boundary is
# unreliable, so nudge a zero-width insertion point open with a throwaway
# character, bookmark that single character, then shrink the bookmarked
# text back down to nothing - leaving a clean, empty "_GoBack" bookmark
# immediately after "Tomcat 7".
$target.InsertAfter("X")
$marker = $d.Range($target.Start, $target.Start + 1)
$d.Bookmarks.Add("_GoBack", $marker)
$bmRange = $d.Bookmarks("_GoBack").Range
$bmRange.Text = ""
